$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.862.23"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "3.744.02"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("D7").Value = "3.742.94"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.93%  "
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.26"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000248"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("D15").Value = "4.371.69"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").Value = "3.743.75"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "68.889.10"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.28"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.86"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +19.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  +12.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.34"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.73%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.52"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.24%  "
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "3.891.39"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").Value = "3.679.70"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.86"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "436.73"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.98"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.97"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.56"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.56"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("D51").Value = "2.779.86"
$ws.Range("E51").Value = "  +1.22%  "
